# Auto-generated edit script applying cryptos.xlsx value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "28.795.97"
Set-TextCell $ws "E2" "  +2.06%  "
Set-TextCell $ws "D3" "1.575.20"
Set-TextCell $ws "E3" "  -0.85%  "
Set-TextCell $ws "D5" "213.32"
Set-TextCell $ws "E5" "  +0.09%  "
Set-TextCell $ws "E6" "  +0.12%  "
Set-TextCell $ws "E7" "  -0.14%  "
Set-TextCell $ws "D8" "44.67"
Set-TextCell $ws "E8" "  +1.89%  "
Set-TextCell $ws "D9" "24.24"
Set-TextCell $ws "E9" "  +1.29%  "
Set-TextCell $ws "E10" "  -1.27%  "
Set-TextCell $ws "E11" "  -1.03%  "
Set-TextCell $ws "E12" "  -0.10%  "
Set-TextCell $ws "D13" "1.799.96"
Set-TextCell $ws "E13" "  -0.88%  "
Set-TextCell $ws "D14" "1.559.01"
Set-TextCell $ws "E14" "  -1.87%  "
Set-TextCell $ws "D15" "28.762.66"
Set-TextCell $ws "E15" "  +1.81%  "
Set-TextCell $ws "E16" "  -1.41%  "
Set-TextCell $ws "E17" "  -1.53%  "
Set-TextCell $ws "D18" "62.51"
Set-TextCell $ws "E18" "  -1.08%  "
Set-TextCell $ws "D19" "230.51"
Set-TextCell $ws "E19" "  +1.57%  "
Set-TextCell $ws "E20" "  -0.77%  "
Set-TextCell $ws "D21" "0.0₃0695"
Set-TextCell $ws "E21" "  -1.82%  "
Set-TextCell $ws "E22" "  -0.06%  "
Set-TextCell $ws "E23" "  -4.84%  "
Set-TextCell $ws "E24" "  -1.36%  "
Set-TextCell $ws "E25" "  +8.77%  "
Set-TextCell $ws "D26" "151.97"
Set-TextCell $ws "E26" "  +0.04%  "
Set-TextCell $ws "E27" "  -0.91%  "
Set-TextCell $ws "E28" "  -1.76%  "
Set-TextCell $ws "E29" "  -2.47%  "
Set-TextCell $ws "E31" "  +2.65%  "
Set-TextCell $ws "E32" "  -2.15%  "
Set-TextCell $ws "D33" "3.22"
Set-TextCell $ws "E33" "  -0.48%  "
Set-TextCell $ws "E34" "  -0.97%  "
Set-TextCell $ws "D35" "1.390.11"
Set-TextCell $ws "E35" "  -0.50%  "
Set-TextCell $ws "E36" "  +2.99%  "
Set-TextCell $ws "D37" "1.54"
Set-TextCell $ws "E37" "  -3.16%  "
Set-TextCell $ws "E38" "  +0.50%  "
Set-TextCell $ws "E39" "  +2.86%  "
Set-TextCell $ws "D40" "0.0166"
Set-TextCell $ws "E41" "  -2.36%  "
Set-TextCell $ws "E42" "  +2.25%  "
Set-TextCell $ws "B43" "ARBITRUM"
Set-TextCell $ws "C43" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D43" "0.797"
Set-TextCell $ws "E43" "  -1.90%  "
Set-TextCell $ws "B44" "PaxDollar"
Set-TextCell $ws "C44" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D44" "1.00"
Set-TextCell $ws "E44" "  -0.08%  "
Set-TextCell $ws "D45" "0.0471"
Set-TextCell $ws "E45" "  +2.81%  "
Set-TextCell $ws "D46" "5.52"
Set-TextCell $ws "E46" "  -1.21%  "
Set-TextCell $ws "D47" "0.966"
Set-TextCell $ws "E47" "  -1.45%  "
Set-TextCell $ws "D48" "63.40"
Set-TextCell $ws "E48" "  -1.30%  "
Set-TextCell $ws "D49" "1.711.80"
Set-TextCell $ws "E49" "  -0.74%  "
Set-TextCell $ws "D50" "86.76"
Set-TextCell $ws "E50" "  -0.73%  "
Set-TextCell $ws "E51" "  -0.33%  "
